$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows per repulled data
$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -8
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -1
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = 4
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = -1
